# Refresh the cryptocurrency price/volume snapshot (D = Price, E = Volume(1h))
# for the rows whose figures moved since the last GitHub Actions run.
#
# Several "Price" values (e.g. 575.14, 160.72 ...) look like plain numbers to
# Excel's auto-detection, but the workbook stores every Price/Volume cell as
# text (openpyxl wrote it that way). To keep those cells text-typed instead
# of being silently coerced to numbers, we briefly force the cell to a text
# number format before assigning the value, then restore the default
# ("Normal") cell style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.720.89'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '3.449.83'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.613'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +12.18%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '3.452.26'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('E12').Value = '  +3.25%  '
$ws.Range('D13').Value = '4.048.30'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = '64.854.52'
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('D18').Value = '3.466.97'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('E23').Value = '  +3.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +9.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.57'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('E35').Value = '  +11.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.71'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.20%  '
$ws.Range('E37').Value = '  +4.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0776'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').Value = '2.962.79'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.65%  '
$ws.Range('E47').Value = '  +1.76%  '
$ws.Range('E48').Value = '  +8.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '308.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.22%  '
$ws.Range('E50').Value = '  +4.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.862'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.52%  '
